$d = $word.ActiveDocument

$replacements = @(
    @{old="69×62=4278"; new="61×56=3416"},
    @{old="85×16=1360"; new="77×59=4543"},
    @{old="34×96=3264"; new="99×42=4158"},
    @{old="71×87=6177"; new="20×95=1900"},
    @{old="72×79=5688"; new="90×63=5670"},
    @{old="86×56=4816"; new="29×82=2378"},
    @{old="64×16=1024"; new="21×93=1953"},
    @{old="48×19=912"; new="17×87=1479"},
    @{old="52×73=3796"; new="44×98=4312"},
    @{old="70×45=3150"; new="71×56=3976"},
    @{old="11×92=1012"; new="12×55=660"},
    @{old="23×93=2139"; new="33×58=1914"},
    @{old="93×27=2511"; new="31×91=2821"},
    @{old="70×53=3710"; new="27×92=2484"},
    @{old="97×95=9215"; new="12×58=696"},
    @{old="97×89=8633"; new="82×31=2542"},
    @{old="37×95=3515"; new="32×91=2912"},
    @{old="97×73=7081"; new="49×69=3381"},
    @{old="16×31=496"; new="75×70=5250"},
    @{old="82×77=6314"; new="30×19=570"},
    @{old="33×41=1353"; new="35×76=2660"},
    @{old="29×15=435"; new="99×90=8910"},
    @{old="45×15=675"; new="46×71=3266"},
    @{old="62×47=2914"; new="98×68=6664"},
    @{old="50×38=1900"; new="89×85=7565"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
